$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new shared strings in the order they must be appended to sharedStrings.xml:
# 34: XXXX, 35: BS, 36: test, 37: kljlj
$ws.Range("K7").Value = "XXXX"
$ws.Range("H4").Value = "BS"
$ws.Range("H6").Value = "test"
$ws.Range("H7").Value = "kljlj"

# Update active cell selection to H8
[void]$ws.Range("H8").Select()
